# "minor: set unused and incorrect 'probviol' learn list column to NaN"
#
# Column J ("prob_violation") on the list_study4_try_counter1 sheet held
# stray 0/1 integers left over from an earlier (incorrect/unused) metric.
# Blank it out to the same "NaN" shared-string sentinel already used by
# column K, clear any stray direct cell formatting some of the J cells
# picked up along the way, and point the sheet's active selection at the
# column that was just edited (J) instead of the old K selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("list_study4_try_counter1")

$rng = $ws.Range("J2:J109")
$rng.Style = "Normal"
$rng.Value = "NaN"

$rng.Select()
